$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Five "Antwortzeit ... Check" paragraphs (one per checklist section):
#    drop the <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/>
#    markers that wrapped the whole paragraph. Text/formatting unchanged.
# ---------------------------------------------------------------------------
$antwortzeitXml = @"
<w:p $W><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Antwortzeit</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/><w:t>Check</w:t></w:r></w:p>
"@

$targets = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Antwortzeit`t`t`t`t`tCheck`r") {
        [void]$targets.Add($p)
    }
}
Write-Output "Antwortzeit paragraphs found: $($targets.Count)"
foreach ($p in $targets) {
    $p.Range.InsertXML($antwortzeitXml)
}

# ---------------------------------------------------------------------------
# 2) "Karte mit Marker existiert ... FEHLT" (section 2 only): FEHLT -> Check,
#    drop the red font colour on that run. Everything else stays the same.
# ---------------------------------------------------------------------------
$karteXml = @"
<w:p $W><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Karte mit Marker existiert</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Check</w:t></w:r></w:p>
"@

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Karte mit Marker existiert`t`t`t`tFEHLT`r") {
        $p.Range.InsertXML($karteXml)
        Write-Output "Replaced 'Karte mit Marker ... FEHLT' paragraph"
        break
    }
}

# ---------------------------------------------------------------------------
# 3) "Austattungsmerkmale exisiteren ... FEHLT" (section 5 only):
#    - fix the typo exisiteren -> existieren (and drop its spellStart/spellEnd)
#    - FEHLT -> Check, dropping the red font colour
#    "Austattungsmerkmale" and the following space stay their own (untouched)
#    runs, exactly like the diff shows.
# ---------------------------------------------------------------------------
$austattungXml = @"
<w:p $W><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Austattungsmerkmale</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>existieren</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Check</w:t></w:r></w:p>
"@

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Austattungsmerkmale exisiteren`t`t`tFEHLT`r") {
        $p.Range.InsertXML($austattungXml)
        Write-Output "Replaced 'Austattungsmerkmale exisiteren ... FEHLT' paragraph"
        break
    }
}

Write-Output "Done"
